$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values with new TPM-derived results (see commit: "update scripts wuth new tpm")
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3399353333333333
$ws.Range("H2").Value = 1.019806
$ws.Range("I2").Value = 0.09929991924017606
$ws.Range("J2").Value = 0.09929991924017606
$ws.Range("M2").Value = 8.033114333333334
$ws.Range("N2").Value = 24.099343
$ws.Range("O2").Value = 0.1374088679258946
$ws.Range("P2").Value = 0.1374088679258946
$ws.Range("Q2").Value = 2.730739398606445
$ws.Range("R2").Value = 24.576654587458
$ws.Range("S2").Value = 0.01364468948792535
$ws.Range("T2").Value = 0.01364468948792535
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3399353333333333
$ws.Range("H3").Value = 1.019806
$ws.Range("I3").Value = 0.09929991924017606
$ws.Range("J3").Value = 0.09929991924017606
$ws.Range("O3").Value = 0.6355200716780686
$ws.Range("P3").Value = 0.6355200716780686
$ws.Range("Q3").Value = 12.62975035404867
$ws.Range("R3").Value = 113.667753186438
$ws.Range("S3").Value = 0.06310709179314311
$ws.Range("T3").Value = 0.06310709179314311
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3399353333333333
$ws.Range("H4").Value = 1.019806
$ws.Range("I4").Value = 0.09929991924017606
$ws.Range("J4").Value = 0.09929991924017606
$ws.Range("M4").Value = 13.27489133333333
$ws.Range("N4").Value = 39.824674
$ws.Range("O4").Value = 0.2270710603960369
$ws.Range("P4").Value = 0.2270710603960369
$ws.Range("Q4").Value = 4.512604610360444
$ws.Range("R4").Value = 40.613441493244
$ws.Range("S4").Value = 0.0225481379591076
$ws.Range("T4").Value = 0.0225481379591076
$ws.Range("I5").Value = 0.4094685684206303
$ws.Range("J5").Value = 0.4094685684206303
$ws.Range("M5").Value = 8.033114333333334
$ws.Range("N5").Value = 24.099343
$ws.Range("O5").Value = 0.1374088679258946
$ws.Range("P5").Value = 0.1374088679258946
$ws.Range("Q5").Value = 11.26035107413055
$ws.Range("R5").Value = 101.343159667175
$ws.Range("S5").Value = 0.05626461243791551
$ws.Range("T5").Value = 0.05626461243791553
$ws.Range("I6").Value = 0.4094685684206303
$ws.Range("J6").Value = 0.4094685684206303
$ws.Range("O6").Value = 0.6355200716780686
$ws.Range("P6").Value = 0.6355200716780686
$ws.Range("S6").Value = 0.2602254939525951
$ws.Range("T6").Value = 0.2602254939525951
$ws.Range("I7").Value = 0.4094685684206303
$ws.Range("J7").Value = 0.4094685684206303
$ws.Range("M7").Value = 13.27489133333333
$ws.Range("N7").Value = 39.824674
$ws.Range("O7").Value = 0.2270710603960369
$ws.Range("P7").Value = 0.2270710603960369
$ws.Range("Q7").Value = 18.60796830240555
$ws.Range("R7").Value = 167.47171472165
$ws.Range("S7").Value = 0.0929784620301197
$ws.Range("T7").Value = 0.09297846203011971
$ws.Range("G8").Value = 1.681642333333333
$ws.Range("H8").Value = 5.044927
$ws.Range("I8").Value = 0.4912315123391937
$ws.Range("J8").Value = 0.4912315123391937
$ws.Range("M8").Value = 8.033114333333334
$ws.Range("N8").Value = 24.099343
$ws.Range("O8").Value = 0.1374088679258946
$ws.Range("P8").Value = 0.1374088679258946
$ws.Range("Q8").Value = 13.50882513144011
$ws.Range("R8").Value = 121.579426182961
$ws.Range("S8").Value = 0.06749956600005372
$ws.Range("T8").Value = 0.06749956600005373
$ws.Range("G9").Value = 1.681642333333333
$ws.Range("H9").Value = 5.044927
$ws.Range("I9").Value = 0.4912315123391937
$ws.Range("J9").Value = 0.4912315123391937
$ws.Range("O9").Value = 0.6355200716780686
$ws.Range("P9").Value = 0.6355200716780686
$ws.Range("Q9").Value = 62.47871513248567
$ws.Range("R9").Value = 562.308436192371
$ws.Range("S9").Value = 0.3121874859323304
$ws.Range("T9").Value = 0.3121874859323304
$ws.Range("G10").Value = 1.681642333333333
$ws.Range("H10").Value = 5.044927
$ws.Range("I10").Value = 0.4912315123391937
$ws.Range("J10").Value = 0.4912315123391937
$ws.Range("M10").Value = 13.27489133333333
$ws.Range("N10").Value = 39.824674
$ws.Range("O10").Value = 0.2270710603960369
$ws.Range("P10").Value = 0.2270710603960369
$ws.Range("Q10").Value = 22.32361923653312
$ws.Range("R10").Value = 200.912573128798
$ws.Range("S10").Value = 0.1115444604068096
$ws.Range("T10").Value = 0.1115444604068096
